# Entity balance sheet: remove NPC "capability" combat math (left side B:D,
# driven by a fixed Test Damage value) is kept, and a mirrored "Attack" table
# (right side H:J) is added that shows how many hits the *player* can take
# from each entity, based on a fixed Player HP value and each entity's Attack.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Test Damage value changed 10 -> 8
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 8

# ---------------------------------------------------------------------------
# 2. New "Player HP :" label + value, to the right of the existing table
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = "Player HP :"
$ws.Range("I2").Value = 40

# ---------------------------------------------------------------------------
# 3. Mirror the "Entity Name / HP / # of Hits" header (row 3) into H3:J3,
#    relabeling the middle column "Attack" since it now represents how hard
#    the entity hits the player rather than the player's fixed test damage.
# ---------------------------------------------------------------------------
$ws.Range("B3:D3").Copy() | Out-Null
$ws.Range("H3:J3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H3").Value = "Entity Name"
$ws.Range("I3").Value = "Attack"
$ws.Range("J3").Value = "# of Hits"

# ---------------------------------------------------------------------------
# 4. Mirror the "Marines" group header (row 4) into H4:J4 and merge it,
#    same as B4:D4.
# ---------------------------------------------------------------------------
$ws.Range("B4:D4").Copy() | Out-Null
$ws.Range("H4:J4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = "Marines"
$ws.Range("H4:J4").Merge() | Out-Null

# Marine
$ws.Range("H5").Value = "Marine"
$ws.Range("I5").Value = 2
$ws.Range("J5").Formula = '=ROUNDDOWN($I$2/I5, 0)'

# Marine Captain
$ws.Range("H6").Value = "Marine Captain"
$ws.Range("I6").Value = 6
$ws.Range("J6").Formula = '=ROUNDDOWN($I$2/I6, 0)'

# Morgan has no attack entry; its I7/J7 simply switch from the blank
# left-aligned style to the blank centered (number) style.
$ws.Range("I7:J7").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 5. New "Pirates" group header mirrored onto the (previously blank) row 8,
#    and the existing row 9 header is mirrored across into H9:J9 as a normal
#    (non-merged) "Pirates" attack row.
# ---------------------------------------------------------------------------
$ws.Range("B9:D9").Copy() | Out-Null
$ws.Range("H8:J8").PasteSpecial(-4122) | Out-Null
$ws.Range("H8").Value = "Pirates"
$ws.Range("H8:J8").Merge() | Out-Null

$ws.Range("H9").Value = "Pirates"
$ws.Range("I9").Value = 2
$ws.Range("J9").Formula = '=ROUNDDOWN($I$2/I9, 0)'

# Pirate Captain
$ws.Range("H10").Value = "Pirate Captain"
$ws.Range("I10").Value = 6
$ws.Range("J10").Formula = '=ROUNDDOWN($I$2/I10, 0)'

# Fat Pirate
$ws.Range("H11").Value = "Fat Pirate"
$ws.Range("I11").Value = 6
$ws.Range("J11").Formula = '=ROUNDDOWN($I$2/I11, 0)'

# row 12 (blank, below Fat Pirate) gets centered blank number styling, like I7/J7.
$ws.Range("I12:J12").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 6. New "Animals" group header mirrored onto the (previously blank) row 13.
# ---------------------------------------------------------------------------
$ws.Range("B14:D14").Copy() | Out-Null
$ws.Range("H13:J13").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").Value = "Animals"
$ws.Range("H13:J13").Merge() | Out-Null

# Kung Fu Dugong
$ws.Range("H14").Value = "Kung Fu Dugong"
$ws.Range("I14").Value = 8
$ws.Range("J14").Formula = '=ROUNDDOWN($I$2/I14, 0)'

# Lapahn
$ws.Range("H15").Value = "Lapahn"
$ws.Range("I15").Value = 6
$ws.Range("J15").Formula = '=ROUNDDOWN($I$2/I15, 0)'

# Lapahn Jump (a second, indented attack belonging to Lapahn)
$ws.Range("H15:J15").Copy() | Out-Null
$ws.Range("H16:J16").PasteSpecial(-4122) | Out-Null
$ws.Range("H16").Value = "Lapahn Jump"
$ws.Range("H16").HorizontalAlignment = -4131   # xlLeft
$ws.Range("H16").IndentLevel = 2
$ws.Range("I16").Value = 6
$ws.Range("J16").Formula = '=ROUNDDOWN($I$2/I16, 0)'

# ---------------------------------------------------------------------------
# 7. The old trailing blank H27:J27 cells are removed entirely.
# ---------------------------------------------------------------------------
$ws.Range("H27:J27").ClearContents() | Out-Null
$ws.Range("H27:J27").ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 8. Column widths for the new H:J, N and Q columns.
# ---------------------------------------------------------------------------
$ws.Columns("H").ColumnWidth = 15.6
$ws.Columns("I").ColumnWidth = 9.1
$ws.Columns("J").ColumnWidth = 10.6
$ws.Columns("N").ColumnWidth = 14.8
$ws.Columns("Q").ColumnWidth = 17.1

# ---------------------------------------------------------------------------
# 9. Active selection moved to C10.
# ---------------------------------------------------------------------------
$ws.Range("C10").Select() | Out-Null

Write-Host "done"
